$d = $word.ActiveDocument

$pairs = @(
    @("572×8=", "781×3="),
    @("346×7=", "642×9="),
    @("571×3=", "992×8="),
    @("108×5=", "775×4="),
    @("594×4=", "483×8="),
    @("665×7=", "945×4="),
    @("955×8=", "656×4="),
    @("863×2=", "506×9="),
    @("575×3=", "228×7="),
    @("449×6=", "628×8="),
    @("701×5=", "164×6="),
    @("908×3=", "384×5="),
    @("578×9=", "680×6="),
    @("389×2=", "686×3="),
    @("769×3=", "289×5="),
    @("534×2=", "396×4="),
    @("398×6=", "798×7="),
    @("310×9=", "326×8="),
    @("799×8=", "735×5="),
    @("858×6=", "971×8="),
    @("849×2=", "933×5="),
    @("190×7=", "569×9="),
    @("767×5=", "486×4="),
    @("159×7=", "923×6="),
    @("271×4=", "675×7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
